$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire column N (the "marker_2" column) and clear its
# contents, which removes the header in N1 and the data values in
# N5:N7 while keeping the column's formatting intact (matches
# "Select column N, press Delete").
$col = $ws.Range("N:N")
$col.Select()
$col.ClearContents()
